# Assignment 1 workbook update — add new group member to the Requirements
# sheet's submitter list (name + e-mail), shifting the existing rows down.
#
# The "Requirements" sheet lists one team member per row in columns A (name)
# and B (e-mail) across rows 1-4. This adds a 5th member — תומר חנניה
# (Tomer Hananya) / Tomer.Hananya@e.braude.ac.il — directly beneath them,
# before the requirements table header, matching the existing layout.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the requirements table header (current row 5),
# pushing the header and every requirement row down by one.
$ws.Rows.Item(5).Insert()

$ws.Range("A5").Value = "תומר חנניה"
$ws.Range("B5").Value = "Tomer.Hananya@e.braude.ac.il"

# Re-anchor the frozen pane beneath the (now one-row-lower) header and
# restore the working selection, same as re-freezing panes after the insert.
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A7").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("I12").Select() | Out-Null
